$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-6 ---

# Row 2: COMBUSTION_FIJA / GAS_NATURAL -> DIESEL / 50 / MENSUAL / date 4/1/2022 -> 10/22/2022
$ws.Range("B2").Value = "DIESEL"
$ws.Range("E2").Value = 44856

# Row 3: LOGISTICA_PRODUCTOS_RESIDUOS / CATEGORIA / MATERIA_PRIMA / MENSUAL / date -> 11/22/2022
$ws.Range("E3").Value = 44887

# Row 4: LOGISTICA_PRODUCTOS_RESIDUOS / MEDIO_TRANSPORTE / CAMION_CARGA / MENSUAL / date -> 11/22/2022
$ws.Range("E4").Value = 44887

# Row 5: LOGISTICA_PRODUCTOS_RESIDUOS / DISTANCIA / 80 / MENSUAL / date -> 11/22/2022
$ws.Range("E5").Value = 44887

# Row 6: LOGISTICA_PRODUCTOS_RESIDUOS / PESO / 800 / MENSUAL / date -> 11/22/2022
$ws.Range("E6").Value = 44887

# --- Add new rows 7-11, same layout/style as row 2 (COMBUSTION_FIJA activity, new fuels) ---
# Match number format of existing date column (E2 uses numFmtId 17 "mmm-yy")
$dateFmt = $ws.Range("E2").NumberFormat()
$ws.Range("E7:E11").NumberFormat = $dateFmt

$ws.Range("A7").Value = "COMBUSTION_FIJA"
$ws.Range("B7").Value = "KEROSENE"
$ws.Range("C7").Value = 50
$ws.Range("D7").Value = "MENSUAL"
$ws.Range("E7").Value = 44856

$ws.Range("A8").Value = "COMBUSTION_FIJA"
$ws.Range("B8").Value = "FUEL_OIL"
$ws.Range("C8").Value = 50
$ws.Range("D8").Value = "MENSUAL"
$ws.Range("E8").Value = 44856

$ws.Range("A9").Value = "COMBUSTION_FIJA"
$ws.Range("B9").Value = "NAFTA"
$ws.Range("C9").Value = 50
$ws.Range("D9").Value = "MENSUAL"
$ws.Range("E9").Value = 44856

$ws.Range("A10").Value = "COMBUSTION_FIJA"
$ws.Range("B10").Value = "CARBON_DE_LEÑA"
$ws.Range("C10").Value = 50
$ws.Range("D10").Value = "MENSUAL"
$ws.Range("E10").Value = 44856

$ws.Range("A11").Value = "COMBUSTION_FIJA"
$ws.Range("B11").Value = "LEÑA"
$ws.Range("C11").Value = 50
$ws.Range("D11").Value = "MENSUAL"
$ws.Range("E11").Value = 44856

# --- Column widths (best fit) matching target cols element ---
# (engine quantizes stored width to the grid it uses internally; offset empirically
# calibrated so the exported <col> widths land as close as possible to target)
$ws.Columns.Item(1).ColumnWidth = 31.592447916666668
$ws.Columns.Item(2).ColumnWidth = 19.022135416666668
$ws.Columns.Item(3).ColumnWidth = 15.166666666666666
$ws.Columns.Item(4).ColumnWidth = 8.736979166666666
$ws.Columns.Item(5).ColumnWidth = 8.451822916666666

# --- Selection moved to G12 like target file ---
$ws.Range("G12").Select()
